$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$xlPasteFormats = -4122

# ---------------------------------------------------------------
# 1. Extend the table with 3 new columns (N, O, P)
# ---------------------------------------------------------------
$lo.ListColumns.Add() | Out-Null
$lo.ListColumns.Add() | Out-Null
$lo.ListColumns.Add() | Out-Null

# Header N1 (new shared string #1)
$ws.Range("N1").Value = "Guarantee A/c No. / #CR"

# Row3 new cell, plain text (new shared string #2)
$ws.Range("N3").Value = "2050088104-CR"

# ---------------------------------------------------------------
# 2. Add a new data row (row 4) to the table, values set in the
#    order they were first typed so the shared-string table lines
#    up with the target workbook.
# ---------------------------------------------------------------
$lo.ListRows.Add() | Out-Null

$ws.Range("F4").Value = "PTP"                       # new shared string #3
$ws.Range("G4").Value = "PTP"                       # reuse #3
$ws.Range("H4").Value = "Email"                     # new shared string #4
$ws.Range("I4").Value = "SME"                       # reuse (existing)
$ws.Range("J4").Value = "Promise to Pay"             # new shared string #5

# K4 needs to be stored as text "2" (quote-prefixed). Apply the
# *source* formatting first (style without quotePrefix) so that the
# apostrophe-driven text entry creates exactly one new cell style.
$ws.Range("K3").Copy() | Out-Null
$ws.Range("K4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("K4").Formula = "'2"                       # new shared string #6 (+ new style)

$ws.Range("A4").Value = "Valid PTP Test 3"           # new shared string #7

# Header O1 (new shared string #8)
$ws.Range("O1").Value = "InstallmentAmount"

$ws.Range("P4").Value = "ACCOUNT TRANSFER"           # new shared string #9

# Header P1 (new shared string #10)
$ws.Range("P1").Value = "PaymentMode"

# Remaining cells reuse already-existing shared strings
$ws.Range("B4").Value = "CO2"
$ws.Range("C4").Value = "<blank>"
$ws.Range("D4").Value = 45547.631249999999
$ws.Range("E4").Value = 45547.631249999999
$ws.Range("L4").Value = "Pre DPD 1- 29 and OS AMT =<3000000"

# M4 and O4 must also be stored as text ("2050088104" / "2"). Apply the
# final target style (which already carries quotePrefix) first, then
# assign via a leading apostrophe so no extra style gets created.
$ws.Range("M3").Copy() | Out-Null
$ws.Range("M4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("M4").Formula = "'2050088104"

$ws.Range("N4").Value = "2050088104-CR"

$ws.Range("M3").Copy() | Out-Null
$ws.Range("O4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("O4").Formula = "'2"

# ---------------------------------------------------------------
# 3. Formatting for all the other (plain text / blank) cells.
#    Values were already assigned above, so pasting formats now is
#    safe and won't create extra cell styles.
# ---------------------------------------------------------------

# Header row N1:P1 - match existing header style (M1)
$ws.Range("M1").Copy() | Out-Null
$ws.Range("N1:P1").PasteSpecial($xlPasteFormats) | Out-Null

# Blank row2 cells in the new columns, formatted like M2
$ws.Range("M2").Copy() | Out-Null
$ws.Range("N2:P2").PasteSpecial($xlPasteFormats) | Out-Null

# Row3 N3 formatted like M3
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial($xlPasteFormats) | Out-Null

# Row3 O3/P3 blank cells, formatted like the other row3 blank style (M2)
$ws.Range("M2").Copy() | Out-Null
$ws.Range("O3:P3").PasteSpecial($xlPasteFormats) | Out-Null

# Row4 formatting, taken from the row above (row 3) which has the
# equivalent look for each column (K4, M4 and O4 were already set above)
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B3:C3").Copy() | Out-Null
$ws.Range("B4:C4").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F3:J3").Copy() | Out-Null
$ws.Range("F4:J4").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("L3:N3").Copy() | Out-Null
$ws.Range("L4:N4").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("L3").Copy() | Out-Null
$ws.Range("P4").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------
# 4. Page setup (printer settings) for the worksheet
# ---------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# 5. Selection
# ---------------------------------------------------------------
$ws.Range("G4").Select() | Out-Null
